$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C32: the trailing non-breaking space (U+00A0) before the line break is
# normalized to a regular space (U+0020).
$ws.Range("C32").Value = '[name=""]   Clear \ Visibility: 14 km ' + [char]10

# C71: wrap the quoted phrase with single quotes instead of double quotes
# (the apostrophe inside "Knights'" stays as-is).
$ws.Range("C71").Value = '[name="Carol"]   ...All because of those ''Knights'' Treasures.''' + [char]10

# C79: replace the surrounding double quotes with single quotes.
$ws.Range("C79").Value = '[name="Carol"]   ''The ancient Knights of Kazimierz were all buried with their riches near the lands they once called home. Their souls protect these unmarked graves for all of eternity.''' + [char]10

# C80: replace the surrounding double quotes with single quotes.
$ws.Range("C80").Value = '[name="Carol"]   ''Only those who fear not sacrifice and possess the true, dauntless bloodline of Kazimierz may open the path.''' + [char]10

# Re-run AutoFit on the touched rows so the COM runtime does not stamp an
# explicit custom row height (the source workbook keeps default row heights).
$ws.Rows(32).AutoFit() | Out-Null
$ws.Rows(71).AutoFit() | Out-Null
$ws.Rows(79).AutoFit() | Out-Null
$ws.Rows(80).AutoFit() | Out-Null
